# Roster update: End-time shift changeover.
#  - "arrival morning" / "departure morning": 21:45 -> 22:00
#  - "arrival night"   / "departure night"  : 09:45 -> 10:00
# (all occurrences live only in column D, the "End" column)
# Also nudges the active selection on "arrival night" to H161.

$wb = $excel.ActiveWorkbook

$wsArrMorning = $wb.Worksheets.Item("arrival morning")
$wsArrMorning.Activate()
$wsArrMorning.Range("D1:D201").Replace("21:45", "22:00", 1, $false, $true) | Out-Null
$excel.ActiveWindow.ScrollRow = 173

$wsArrNight = $wb.Worksheets.Item("arrival night")
$wsArrNight.Activate()
$wsArrNight.Range("D1:D161").Replace("09:45", "10:00", 1, $false, $true) | Out-Null
$wsArrNight.Range("H161").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 131

$wsDepMorning = $wb.Worksheets.Item("departure morning")
$wsDepMorning.Activate()
$wsDepMorning.Range("D1:D181").Replace("21:45", "22:00", 1, $false, $true) | Out-Null
$excel.ActiveWindow.ScrollRow = 167

$wsDepNight = $wb.Worksheets.Item("departure night")
$wsDepNight.Activate()
$wsDepNight.Range("D1:D129").Replace("09:45", "10:00", 1, $false, $true) | Out-Null
